$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 (I0) and J1 (IF), matching the header style used by the
# existing columns (e.g. H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for I and J columns for data rows 2-37: each entry is (I, J).
$ijValues = @(
    @(1, 7),
    @(1, 7),
    @(1, 4),
    @(1, 7),
    @(3, 9),
    @(1, 7),
    @(1, 4),
    @(1, 7),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 7),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(3, 8),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(1, 2),
    @(3, 3)
)

$startRow = 2
for ($idx = 0; $idx -lt $ijValues.Count; $idx++) {
    $row = $startRow + $idx
    $pair = $ijValues[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
